# إضافة حدث جديد في Card4
# Adds a new event row (row 22) to the "Card4" worksheet, mirroring the
# existing event-log rows (columns A, L, M, N, O filled; B:K left blank).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card4")

# Column A holds the card number as text ("4") in every existing row, so
# force text storage (leading apostrophe) instead of letting it become a
# numeric 4, then strip the resulting style tweak back to Normal.
$cardCell = $ws.Cells.Item(22, 1)
$cardCell.Value = "'4"
$cardCell.Style = "Normal"

$ws.Cells.Item(22, 12).Value = "21/1/2026"
$ws.Cells.Item(22, 13).Value = "قطع سير700"
$ws.Cells.Item(22, 14).Value = "تم تغير سير700 (مشلان)"
$ws.Cells.Item(22, 15).Value = "محمود ايهاب"
